$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-20 22:48:40'
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '63%'
$ws.Range("N2").Value = '-2.9 °C 22:09 TU'
$ws.Range("E3").Value = '2026-02-20 22:48:43'
$ws.Range("E4").Value = '2026-02-20 22:48:46'
$ws.Range("J4").Value = '1023.1 hPa'
$ws.Range("N4").Value = '3.8 °C 22:29 TU'
$ws.Range("O4").Value = '9.7 °C'
$ws.Range("E5").Value = '2026-02-20 22:48:49'
$ws.Range("E6").Value = '2026-02-20 22:48:51'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '69%'
$ws.Range("J6").Value = '1023.0 hPa'
$ws.Range("L6").Value = '16.9 km/h - 331º 22:16 TU'
$ws.Range("O6").Value = '9.3 °C'
$ws.Range("E7").Value = '2026-02-20 22:48:54'
$ws.Range("J7").Value = '1022.9 hPa'
$ws.Range("K7").Value = '12.0 MJ/m2'
$ws.Range("O7").Value = '13.3 °C'
$ws.Range("E8").Value = '2026-02-20 22:48:57'
$ws.Range("J8").Value = '1023.2 hPa'
$ws.Range("E9").Value = '2026-02-20 22:49:00'
$ws.Range("E10").Value = '2026-02-20 22:49:02'
$ws.Range("O10").Value = '7.3 °C'
$ws.Range("E11").Value = '2026-02-20 22:49:05'
$ws.Range("E12").Value = '2026-02-20 22:49:08'
$ws.Range("E13").Value = '2026-02-20 22:49:10'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '47%'
$ws.Range("J13").Value = '1024.4 hPa'
$ws.Range("N13").Value = '-1.1 °C 22:26 TU'
$ws.Range("O13").Value = '5.8 °C'
$ws.Range("E14").Value = '2026-02-20 22:49:13'
$ws.Range("O14").Value = '11.7 °C'
$ws.Range("E15").Value = '2026-02-20 22:49:15'
$ws.Range("E16").Value = '2026-02-20 22:49:18'
$ws.Range("E17").Value = '2026-02-20 22:49:21'
$ws.Range("O17").Value = '3.3 °C'
$ws.Range("E18").Value = '2026-02-20 22:49:23'
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = '79%'
$ws.Range("J18").Value = '1023.3 hPa'
$ws.Range("O18").Value = '7.6 °C'
$ws.Range("E19").Value = '2026-02-20 22:49:26'
$ws.Range("E20").Value = '2026-02-20 22:49:28'
$ws.Range("E21").Value = '2026-02-20 22:49:31'
$ws.Range("J21").Value = '1023.3 hPa'
$ws.Range("N21").Value = '2.8 °C 22:25 TU'
$ws.Range("O21").Value = '8.9 °C'
$ws.Range("E22").Value = '2026-02-20 22:49:34'
$ws.Range("O22").Value = '-3.9 °C'
$ws.Range("E23").Value = '2026-02-20 22:49:37'
$ws.Range("O23").Value = '-4.5 °C'
$ws.Range("E24").Value = '2026-02-20 22:49:39'
$ws.Range("J24").Value = '1025.7 hPa'
$ws.Range("E25").Value = '2026-02-20 22:49:42'
$ws.Range("O25").Value = '-1.3 °C'
$ws.Range("E26").Value = '2026-02-20 22:49:44'
$ws.Range("O26").Value = '5.6 °C'
$ws.Range("E27").Value = '2026-02-20 22:49:47'
$ws.Range("E28").Value = '2026-02-20 22:49:50'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '67%'
$ws.Range("J28").Value = '1023.5 hPa'
$ws.Range("O28").Value = '6.8 °C'
$ws.Range("E29").Value = '2026-02-20 22:49:52'
$ws.Range("E30").Value = '2026-02-20 22:49:55'
$ws.Range("J30").Value = '1022.8 hPa'
$ws.Range("E31").Value = '2026-02-20 22:49:58'
$ws.Range("J31").Value = '1022.0 hPa'
$ws.Range("E32").Value = '2026-02-20 22:50:01'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '84%'
$ws.Range("O32").Value = '4.4 °C'
$ws.Range("E33").Value = '2026-02-20 22:50:03'
$ws.Range("J33").Value = '1023.6 hPa'
$ws.Range("N33").Value = '1.6 °C 22:24 TU'
$ws.Range("E34").Value = '2026-02-20 22:50:06'
$ws.Range("O34").Value = '1.2 °C'
$ws.Range("E35").Value = '2026-02-20 22:50:08'
$ws.Range("J35").Value = '1027.1 hPa'
$ws.Range("E36").Value = '2026-02-20 22:50:11'
$ws.Range("J36").Value = '1022.9 hPa'
$ws.Range("E37").Value = '2026-02-20 22:50:14'
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = '69%'
$ws.Range("J37").Value = '1025.1 hPa'
$ws.Range("O37").Value = '4.4 °C'
$ws.Range("E38").Value = '2026-02-20 22:50:17'
$ws.Range("E39").Value = '2026-02-20 22:50:19'
$ws.Range("K39").Value = '14.7 MJ/m2'
$ws.Range("E40").Value = '2026-02-20 22:50:22'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '39%'
$ws.Range("J40").Value = '1024.0 hPa'
$ws.Range("O40").Value = '10.0 °C'
$ws.Range("E41").Value = '2026-02-20 22:50:25'
$ws.Range("J41").Value = '1023.4 hPa'
$ws.Range("O41").Value = '13.3 °C'
$ws.Range("E42").Value = '2026-02-20 22:50:28'
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = '64%'
$ws.Range("E43").Value = '2026-02-20 22:50:30'
$ws.Range("E44").Value = '2026-02-20 22:50:33'
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = '74%'
$ws.Range("M44").Value = '-0.3 °C 22:29 TU'
$ws.Range("O44").Value = '-4.3 °C'
$ws.Range("E45").Value = '2026-02-20 22:50:36'
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '84%'
$ws.Range("J45").Value = '1030.0 hPa'
$ws.Range("E46").Value = '2026-02-20 22:50:38'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '57%'
$ws.Range("J46").Value = '1026.7 hPa'
$ws.Range("N46").Value = '8.2 °C 22:28 TU'
$ws.Range("O46").Value = '11.8 °C'
